$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 173, shifting existing rows 173..184 down to 174..185.
$ws.Rows.Item(173).Insert()

# Update row 172: new record date (D) and unit of sale (N) change; other columns unchanged.
$ws.Range("D172").Value = 44585
$ws.Range("N172").Value = "$/caja 10 kilos"

# Populate the freshly inserted row 173 with the record that used to live in row 172
# (same data values, previous date and unit-of-sale).
$ws.Range("A173").Value = 7
$ws.Range("B173").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C173").Value = "Ñuble"
$ws.Range("D173").Value = 44560
$ws.Range("E173").Value = 16
$ws.Range("F173").Value = 100112003
$ws.Range("G173").Value = "Ajo"
$ws.Range("H173").Value = "Chino"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 60
$ws.Range("K173").Value = 19000
$ws.Range("L173").Value = 20000
$ws.Range("M173").Value = 19500
$ws.Range("N173").Value = "$/malla 10 kilos"
$ws.Range("O173").Value = "China"
$ws.Range("P173").Value = 1950
$ws.Range("Q173").Value = 10
$ws.Range("R173").Value = "Hortaliza"
